$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.864.63"
$ws.Range("E2").Value = "  +1.87%  "
$ws.Range("D3").Value = "3.034.53"
$ws.Range("E3").Value = "  +0.81%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.37"
$ws.Range("E5").Value = "  -0.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.98"
$ws.Range("E6").Value = "  +6.13%  "
$ws.Range("D8").Value = "3.033.06"
$ws.Range("E8").Value = "  +0.77%  "
$ws.Range("E9").Value = "  -1.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.71"
$ws.Range("E10").Value = "  +12.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.152"
$ws.Range("E11").Value = "  +2.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.465"
$ws.Range("E12").Value = "  +1.11%  "
$ws.Range("E13").Value = "  +1.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.64"
$ws.Range("E14").Value = "  +3.14%  "
$ws.Range("D16").Value = "3.535.05"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.07"
$ws.Range("E17").Value = "  +0.42%  "
$ws.Range("D18").Value = "62.824.18"
$ws.Range("E18").Value = "  +1.94%  "
$ws.Range("D19").Value = "3.031.75"
$ws.Range("E19").Value = "  +0.87%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "453.73"
$ws.Range("E20").Value = "  -0.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.29"
$ws.Range("E21").Value = "  +1.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.698"
$ws.Range("E22").Value = "  +1.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.50"
$ws.Range("E23").Value = "  +1.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.15"
$ws.Range("E24").Value = "  +1.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.16"
$ws.Range("E25").Value = "  +3.36%  "
$ws.Range("E26").Value = "  +2.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.29"
$ws.Range("E27").Value = "  +2.13%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.45"
$ws.Range("E29").Value = "  +2.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.70"
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.24"
$ws.Range("E31").Value = "  +7.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.63"
$ws.Range("E33").Value = "  +0.19%  "
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("D35").Value = "0.0₃0868"
$ws.Range("E35").Value = "  +2.40%  "
$ws.Range("E36").Value = "  +1.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.93"
$ws.Range("E37").Value = "  +2.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.20"
$ws.Range("E38").Value = "  +9.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.10"
$ws.Range("E39").Value = "  +0.82%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.129"
$ws.Range("E40").Value = "  +3.22%  "
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.45"
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("E42").Value = "  -1.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.302"
$ws.Range("E43").Value = "  +10.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.80"
$ws.Range("E44").Value = "  +5.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "394.03"
$ws.Range("E45").Value = "  -2.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0360"
$ws.Range("E46").Value = "  +1.29%  "
$ws.Range("D47").Value = "2.727.20"
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "131.94"
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.31"
$ws.Range("E49").Value = "  +6.19%  "
$ws.Range("B50").Value = "USDe"
$ws.Range("C50").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.47"
$ws.Range("E51").Value = "  +3.52%  "
